$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A3").Value = "130, 426, 458, SF"

$ws.Range("A8").Value = "98, 130, 754, 786"

$ws.Range("A9").Value = "130, 754, 1082, 1114"

$ws.Range("A10").Value = "426, 754, 1082, 1082"
$ws.Range("C10").Value = "754, 1082, 1082, 426"

$ws.Range("A12").Value = "98, 130, 426, 426, 1082"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "426"

$ws.Range("A15").Value = "98, 426, 426, 754, 786"

$ws.Range("A16").Value = "98, 130, 130, 426, 754, 1082"
